# This script applies the edits described by the diff:
#  1. Remove the "Bentleigh" row (all other rows shift up by one).
#  2. Update the "Chadstone" row's Notes text.
#  3. Insert a new "Hallam" row (just before "Lakes Entrance").
#  4. Update the "Moorabbin" row's Exposure period text.
#  5. Insert a new "Mordialloc" row (just before "Mount Martha").
#  6. Append a new "Wonthaggi" row at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the "Bentleigh" row entirely.
$bentleigh = $ws.Columns.Item(1).Find("Bentleigh")
$ws.Rows.Item($bentleigh.Row).Delete()

# 2. Update the "Chadstone" row's Notes (column D).
$chadstone = $ws.Columns.Item(1).Find("Chadstone")
$ws.Cells.Item($chadstone.Row, 4).Value = "Case did not attend during infectious period but may have acquired their illness here. If you attended Chadstone Shopping Centre but did not attend to any of the acquisition site stores listed above, monitor for symptoms - If symptoms develop, immediately get tested and isolate until you receive a negative result."

# 3. Insert a new "Hallam" row just before "Lakes Entrance".
$lakesEntrance = $ws.Columns.Item(1).Find("Lakes Entrance")
$newRow = $lakesEntrance.Row
$ws.Rows.Item($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "Hallam"
$ws.Cells.Item($newRow, 2).Value = "Coles Hallam, 2 Princes Domain Drive, Hallam, VIC 3803"
$ws.Cells.Item($newRow, 3).Value = "30/12/20 6:15am - 6:30am"
$ws.Cells.Item($newRow, 4).Value = "Case shopped in store"

# 4. Update the "Moorabbin" row's Exposure period (column C).
$moorabbin = $ws.Columns.Item(1).Find("Moorabbin")
$ws.Cells.Item($moorabbin.Row, 3).Value = "30/12/20 10:45am - 12:15pm and 4:00pm- 5:50pm"

# 5. Insert a new "Mordialloc" row just before "Mount Martha".
$mountMartha = $ws.Columns.Item(1).Find("Mount Martha")
$newRow2 = $mountMartha.Row
$ws.Rows.Item($newRow2).Insert()
$ws.Cells.Item($newRow2, 1).Value = "Mordialloc"
$ws.Cells.Item($newRow2, 2).Value = "Woodlands Golf Club, 109 White Street, Mordialloc, VIC 3195"
$ws.Cells.Item($newRow2, 3).Value = "28/12/20 12:00pm - 6:00pm"
$ws.Cells.Item($newRow2, 4).Value = "Case attended course"

# 6. Append a new "Wonthaggi" row after "Southern Cross" (the last row).
$southernCross = $ws.Columns.Item(1).Find("Southern Cross")
$newRow3 = $southernCross.Row + 1
$ws.Cells.Item($newRow3, 1).Value = "Wonthaggi"
$ws.Cells.Item($newRow3, 2).Value = "Wonthaggi Plaza Shopping Centre, 2 Biggs Drive, Wonthaggi, VIC 3995"
$ws.Cells.Item($newRow3, 3).Value = "28/12/20 1:30pm - 2:30pm"
$ws.Cells.Item($newRow3, 4).Value = "Kmart - shopped for 15 mins"

# Update the sheet dimension to reflect the new extent of the data.
$ws.UsedRange | Out-Null
